$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy formatting from column F (rows 11-16) into new column G
$ws.Range("F11:F16").Copy()
$ws.Range("G11:G16").PasteSpecial(-4122)

# Set the new "Login" column values
$ws.Cells.Item(12, 7).Value = "Login"
$ws.Cells.Item(13, 7).Value = "Login"
$ws.Cells.Item(14, 7).Value = "soufian"
$ws.Cells.Item(15, 7).Value = "mathis"
$ws.Cells.Item(16, 7).Value = "dimitri"

$ws.Range("G16").Select()
